# Generate Report for Handback
#
# The "96e78804-1ed0-4da2-b7d2-3f6c213bc972.md" file has now been handed
# back (it is in sync with en-US), so flip its status away from
# "Ready for handoff" everywhere it is reported, and stamp the per-locale
# handback datetime for the locales that just completed.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# Overview sheet: columns B (zh-cn) and C (de-de) hold the status for this file.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status

# zh-cn detail sheet: column C is Status, column H is Latest Handback DateTime.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("H3").Value = "2016-03-13 08:51:26"

# de-de detail sheet: column C is Status, column H is Latest Handback DateTime.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("H3").Value = "2016-03-13 08:51:32"
